# Add a new tour ("file_name") pair of columns (AH label / AI cost) to Sheet2,
# mirroring the existing per-tour column pairs (C/D, G/H, ..., AE/AF).
# This corresponds to upstream commit "made file_name user changeable":
# a new tour's step costs + a running total are appended in columns AH:AI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "0  . Harrisburg          -->       1  . Charleston          Cost",
    "1  . Charleston          -->       2  . Columbus            Cost",
    "2  . Columbus            -->       3  . Frankfort           Cost",
    "3  . Frankfort           -->       4  . Austin              Cost",
    "4  . Austin              -->       5  . Santa Fe            Cost",
    "5  . Santa Fe            -->       6  . Salt Lake City      Cost",
    "6  . Salt Lake City      -->       7  . Boise               Cost",
    "7  . Boise               -->       8  . Carson City         Cost",
    "8  . Carson City         -->       9  . Phoenix             Cost",
    "9  . Phoenix             -->       10 . Albany              Cost",
    "10 . Albany              -->       11 . Boston              Cost",
    "11 . Boston              -->       12 . Providence          Cost",
    "12 . Providence          -->       13 . Columbia            Cost",
    "13 . Columbia            -->       14 . Tallahassee         Cost",
    "14 . Tallahassee         -->       15 . Montgomery          Cost",
    "15 . Montgomery          -->       16 . Little Rock         Cost",
    "16 . Little Rock         -->       17 . Baton Rouge         Cost",
    "17 . Baton Rouge         -->       18 . Jackson             Cost",
    "18 . Jackson             -->       19 . Oklahoma City       Cost",
    "19 . Oklahoma City       -->       20 . Lincoln             Cost",
    "20 . Lincoln             -->       21 . Denver              Cost",
    "21 . Denver              -->       22 . Cheyenne            Cost",
    "22 . Cheyenne            -->       23 . Olympia             Cost",
    "23 . Olympia             -->       24 . Juneau              Cost",
    "24 . Juneau              -->       25 . Honolulu            Cost",
    "25 . Honolulu            -->       26 . Sacramento          Cost",
    "26 . Sacramento          -->       27 . Salem               Cost",
    "27 . Salem               -->       28 . Helana              Cost",
    "28 . Helana              -->       29 . Bismarck            Cost",
    "29 . Bismarck            -->       30 . Pierre              Cost",
    "30 . Pierre              -->       31 . Saint Paul          Cost",
    "31 . Saint Paul          -->       32 . Lansing             Cost",
    "32 . Lansing             -->       33 . Montpelier          Cost",
    "33 . Montpelier          -->       34 . Augusta             Cost",
    "34 . Augusta             -->       35 . Concord             Cost",
    "35 . Concord             -->       36 . Hartford            Cost",
    "36 . Hartford            -->       37 . Trenton             Cost",
    "37 . Trenton             -->       38 . Indianapolis        Cost",
    "38 . Indianapolis        -->       39 . Springfield         Cost",
    "39 . Springfield         -->       40 . Madison             Cost",
    "40 . Madison             -->       41 . Des Moines          Cost",
    "41 . Des Moines          -->       42 . Topeka              Cost",
    "42 . Topeka              -->       43 . Jefferson City      Cost",
    "43 . Jefferson City      -->       44 . Nashville           Cost",
    "44 . Nashville           -->       45 . Atlanta             Cost",
    "45 . Atlanta             -->       46 . Raleigh             Cost",
    "46 . Raleigh             -->       47 . Richmond            Cost",
    "47 . Richmond            -->       48 . Annapolis           Cost",
    "48 . Annapolis           -->       49 . Dover               Cost",
    "49 . Dover               -->       0  . Harrisburg          Cost"
)

$values = @(5.13, 2.11, 2.57, 15.13, 9.83, 7.81, 5.2, 5.67, 9.57, 39.380000000000003, 2.79, 0.56999999999999995, 12.39, 4.8, 2.77, 6.5, 4.4400000000000004, 2.08, 7.98, 5.4, 8.3800000000000008, 1.42, 19.03, 16.11, 43.78, 40.24, 6.56, 11.13, 11.46, 4.47, 7.27, 8.83, 12.07, 2.81, 2.1, 1.84, 2.59, 11.4, 3.5, 3.3, 4.49, 3.28, 3.53, 5.92, 3.39, 6.09, 2.13, 1.72, 0.99, 1.75)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 34).Value = $labels[$i]
    $ws.Cells.Item($row, 35).Value = $values[$i]
}

# Running total for the new tour, styled like the existing Z53 total (2 d.p.).
$ws.Range("AI53").Formula = "=SUM(AI3:AI52)"
$ws.Range("AI53").NumberFormat = "0.00"

# Match the author's final selection/scroll position.
$ws.Range("P28").Select() | Out-Null
$ws.Range("AL44").Select() | Out-Null
